$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure values are written as literal text (matching original inline string cells)
# even when the text looks like a number (e.g. "599.90" or "0.999").
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "67.755.97"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "2.675.09"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "599.90"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "167.17"
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "2.674.99"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "0.144"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "0.363"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "5.23"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").Value = "27.91"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "3.164.23"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "0.0000185"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").Value = "67.625.92"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "2.689.57"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "11.76"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "7.79"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "364.36"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("E22").Value = "  -3.25%  "
$ws.Range("D23").Value = "4.84"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -3.66%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "70.95"
$ws.Range("E26").Value = "  -4.05%  "
$ws.Range("D27").Value = "10.21"
$ws.Range("E27").Value = "  +3.03%  "
$ws.Range("D28").Value = "2.810.08"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "557.52"
$ws.Range("E31").Value = "  -6.36%  "
$ws.Range("D32").Value = "8.04"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").Value = "1.40"
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "1.56"
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("D38").Value = "19.56"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").Value = "0.373"
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("D41").Value = "5.32"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("E42").Value = "  -4.13%  "
$ws.Range("D43").Value = "17.95"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "2.53"
$ws.Range("E44").Value = "  -6.24%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "40.31"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("E47").Value = "  -4.65%  "
$ws.Range("D48").Value = "0.592"
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("D49").Value = "153.76"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("D50").Value = "3.87"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("E51").Value = "  -3.21%  "

# Restore default (General) style so cells keep looking like the originals
$rng.Style = "Normal"

